$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Choose Flight" value for the bookFlight row now carries an "Airline:" prefix.
$ws.Range("E2").Value = "Airline:Virgin America"

# Column E was auto best-fit by Excel to accommodate the longer text.
$ws.Range("E1").ColumnWidth = 20

# The saved cursor/selection moved to E3 (the newly edited column).
$ws.Range("E3").Select()
